$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet references
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# Grow each table by two rows (this also grows dimension/used range).
# ---------------------------------------------------------------------
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null
$loOverview.ListRows.Add() | Out-Null

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null
$loZhCn.ListRows.Add() | Out-Null

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null
$loDeDe.ListRows.Add() | Out-Null

# ---------------------------------------------------------------------
# Overview sheet - two new handed-off files
# ---------------------------------------------------------------------
$wsOverview.Range("A6").Value = "50af4f54-a5ca-4fe5-b227-605c9730f8a3.md"
$wsOverview.Range("C6").Value = ".md"
$wsOverview.Range("D6").Value = "'"
$wsOverview.Range("E6").Value = "Ready for handoff"
$wsOverview.Range("F6").Value = "Ready for handoff"
$wsOverview.Range("G6").Value = "2016-08-12 04:43:51"

$wsOverview.Range("A7").Value = "a0229e38-a74f-4ea5-a35d-8ba91c986ce5.md"
$wsOverview.Range("C7").Value = ".md"
$wsOverview.Range("D7").Value = "'"
$wsOverview.Range("E7").Value = "Ready for handoff"
$wsOverview.Range("F7").Value = "Ready for handoff"
$wsOverview.Range("G7").Value = "2016-08-12 04:43:51"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B6"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/50af4f54a5ca4fe5b227605c9730f8a3000000/e2e/50af4f54-a5ca-4fe5-b227-605c9730f8a3.md", "", "", "e2e\50af4f54-a5ca-4fe5-b227-605c9730f8a3.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B7"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/a0229e38a74f4ea5a35d8ba91c986ce5000000/e2e/a0229e38-a74f-4ea5-a35d-8ba91c986ce5.md", "", "", "e2e\a0229e38-a74f-4ea5-a35d-8ba91c986ce5.md") | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet - detail rows for the same two files
# ---------------------------------------------------------------------
$wsZhCn.Range("B6").Value = ".md"
$wsZhCn.Range("C6").Value = "Ready for handoff"
$wsZhCn.Range("D6").Value = "e2e"
$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("F6").Value = "'False"
$wsZhCn.Range("G6").Value = "50af4f54-a5ca-4fe5-b227-605c9730f8a3.7bef278b6e0e5db80a9f61f475abc0e326c9271c.zh-cn.xlf"
$wsZhCn.Range("H6").Value = "2016-08-12 04:43:46"
$wsZhCn.Range("I6").Value = "'"
$wsZhCn.Range("J6").Value = "'"
$wsZhCn.Range("K6").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L6").Value = "'"
$wsZhCn.Range("M6").Value = "'True"
$wsZhCn.Range("N6").Value = "'"
$wsZhCn.Range("O6").Value = "'False"
$wsZhCn.Range("P6").Value = "'"

$wsZhCn.Range("B7").Value = ".md"
$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Range("D7").Value = "e2e"
$wsZhCn.Range("E7").Value = "ht"
$wsZhCn.Range("F7").Value = "'False"
$wsZhCn.Range("G7").Value = "a0229e38-a74f-4ea5-a35d-8ba91c986ce5.5dcefbe005e1319249a2ffd77168c30d2d53363a.zh-cn.xlf"
$wsZhCn.Range("H7").Value = "2016-08-12 04:43:46"
$wsZhCn.Range("I7").Value = "'"
$wsZhCn.Range("J7").Value = "'"
$wsZhCn.Range("K7").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L7").Value = "'"
$wsZhCn.Range("M7").Value = "'True"
$wsZhCn.Range("N7").Value = "'"
$wsZhCn.Range("O7").Value = "'False"
$wsZhCn.Range("P7").Value = "'"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A6"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/50af4f54a5ca4fe5b227605c9730f8a3000000/e2e/50af4f54-a5ca-4fe5-b227-605c9730f8a3.md", "", "", "50af4f54-a5ca-4fe5-b227-605c9730f8a3.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A7"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/a0229e38a74f4ea5a35d8ba91c986ce5000000/e2e/a0229e38-a74f-4ea5-a35d-8ba91c986ce5.md", "", "", "a0229e38-a74f-4ea5-a35d-8ba91c986ce5.md") | Out-Null

# ---------------------------------------------------------------------
# de-de sheet - detail rows for the same two files
# ---------------------------------------------------------------------
$wsDeDe.Range("B6").Value = ".md"
$wsDeDe.Range("C6").Value = "Ready for handoff"
$wsDeDe.Range("D6").Value = "e2e"
$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("F6").Value = "'False"
$wsDeDe.Range("G6").Value = "50af4f54-a5ca-4fe5-b227-605c9730f8a3.7bef278b6e0e5db80a9f61f475abc0e326c9271c.de-de.xlf"
$wsDeDe.Range("H6").Value = "2016-08-12 04:43:51"
$wsDeDe.Range("I6").Value = "'"
$wsDeDe.Range("J6").Value = "'"
$wsDeDe.Range("K6").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L6").Value = "'"
$wsDeDe.Range("M6").Value = "'True"
$wsDeDe.Range("N6").Value = "'"
$wsDeDe.Range("O6").Value = "'False"
$wsDeDe.Range("P6").Value = "'"

$wsDeDe.Range("B7").Value = ".md"
$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Range("D7").Value = "e2e"
$wsDeDe.Range("E7").Value = "ht"
$wsDeDe.Range("F7").Value = "'False"
$wsDeDe.Range("G7").Value = "a0229e38-a74f-4ea5-a35d-8ba91c986ce5.5dcefbe005e1319249a2ffd77168c30d2d53363a.de-de.xlf"
$wsDeDe.Range("H7").Value = "2016-08-12 04:43:51"
$wsDeDe.Range("I7").Value = "'"
$wsDeDe.Range("J7").Value = "'"
$wsDeDe.Range("K7").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L7").Value = "'"
$wsDeDe.Range("M7").Value = "'True"
$wsDeDe.Range("N7").Value = "'"
$wsDeDe.Range("O7").Value = "'False"
$wsDeDe.Range("P7").Value = "'"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A6"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/50af4f54a5ca4fe5b227605c9730f8a3000000/e2e/50af4f54-a5ca-4fe5-b227-605c9730f8a3.md", "", "", "50af4f54-a5ca-4fe5-b227-605c9730f8a3.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A7"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/a0229e38a74f4ea5a35d8ba91c986ce5000000/e2e/a0229e38-a74f-4ea5-a35d-8ba91c986ce5.md", "", "", "a0229e38-a74f-4ea5-a35d-8ba91c986ce5.md") | Out-Null
